# Update row 22 (student record #19) on the active sheet:
#  - fill in the previously empty G22/H22 scores with 5
#  - add a new J22 score of 5, matching the formatting already used for
#    the same "J" column on sibling rows (e.g. J9/J12/J17 use style index 9)
#  - L22's SUM formula recalculates automatically from 25 to 40

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("G22").Value = 5
$ws.Range("H22").Value = 5
$ws.Range("J22").Value = 5

# Copy the formatting used by the other "J" column cells (e.g. J9) onto J22
# so the new value picks up the same style (border/fill/alignment) instead
# of the workbook default.
$ws.Range("J9").Copy() | Out-Null
$ws.Range("J22").PasteSpecial(-4122) | Out-Null

# Move the active selection to K11, matching the saved view state.
$ws.Range("K11").Select() | Out-Null
